# Validated UE 2000-2018 against previous work
# Fill in previously-blank ("None") UE values for 7 states/territory rows on
# the DATA sheet, and update the downstream REP/DEM/OTH/TOT/REP_UE/DEM_UE/
# NET_UE summary rows to reflect them. The Analysis sheet just formula-refs
# DATA, so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L")

function Set-RowValues {
    param($Sheet, $Row, $ColNames, $Values)
    for ($i = 0; $i -lt $ColNames.Length; $i++) {
        if ($null -ne $Values[$i]) {
            $Sheet.Range($ColNames[$i] + $Row).Value = $Values[$i]
        }
    }
}

# State rows that previously held "None" placeholders in C:L
Set-RowValues $ws 3  $cols @(0, 0, 0, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5)                 # Alaska
Set-RowValues $ws 9  $cols @(0.5, 0.5, 0.5, 0.5, 0.5, -0.5, -0.5, -0.5, -0.5, -0.5)      # Delaware
Set-RowValues $ws 27 $cols @(0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5)           # Montana
Set-RowValues $ws 35 $cols @(-0.5, -0.5, -0.5, -0.5, -0.5, 0.5, 0.5, 0.5, 0.5, 0.5)      # North Dakota
Set-RowValues $ws 42 $cols @(0.5, 0.5, -0.5, -0.5, -0.5, 0.5, 0.5, 0.5, 0.5, 0.5)        # South Dakota
Set-RowValues $ws 51 $cols @(0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5)           # Wyoming

# Row 46 (territory) only gets F:L filled in; C:E stay "None"
Set-RowValues $ws 46 $cols @($null, $null, $null, -0.5, -0.5, -0.5, 0, -0.5, 0, -0.5)

# Recomputed summary rows (REP, DEM, OTH, TOT, REP_UE, DEM_UE, NET_UE)
Set-RowValues $ws 52 $cols @(221, 229, 233, 202, 178, 242, 234, 247, 241, 199)
Set-RowValues $ws 53 $cols @(212, 205, 201, 233, 257, 193, 201, 188, 194, 235)
Set-RowValues $ws 54 $cols @(2, 1, 1, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues $ws 55 $cols @(435, 435, 435, 435, 435, 435, 435, 435, 435, 435)
Set-RowValues $ws 56 $cols @(26, 27, 31.5, 24, 17, 36.5, 45.5, 44.5, 43.5, 35.5)
Set-RowValues $ws 57 $cols @(-17.5, -22.5, -18, -25.5, -33.5, -22, -25.5, -25, -21.5, -39)
Set-RowValues $ws 58 $cols @(8.5, 4.5, 13.5, -1.5, -16.5, 14.5, 20, 19.5, 22, -3.5)

# Leave the cursor where the author last left it, on the Analysis sheet
$wsAnalysis = $wb.Worksheets.Item("Analysis")
$wsAnalysis.Range("K66").Select()
